# Generate Report for Handback
# Updates the Overview/zh-cn/de-de sheets to reflect that both files have
# been handed back and are back in sync with en-US.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Status column updates (Overview + both locale sheets) ---
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

# --- zh-cn sheet: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3d274448feffc1782fe7b51ae4b7508c91d149bb/e2e/8351b26a-afed-4463-8200-ed7236f08299.md",
    "",
    "",
    "8351b26a-afed-4463-8200-ed7236f08299.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3538982758f4a116054fa89e59de4f03d20fb81b/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/8351b26a-afed-4463-8200-ed7236f08299.1f1eec3eddfab6c04388319c08bf8fbb2f1611cc.zh-cn.xlf",
    "",
    "",
    "8351b26a-afed-4463-8200-ed7236f08299.1f1eec3eddfab6c04388319c08bf8fbb2f1611cc.zh-cn.xlf"
)
$wsZhCn.Range("H2").Value = "2016-03-22 12:11:49"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3d274448feffc1782fe7b51ae4b7508c91d149bb/e2e/c4e3a62b-d9fe-4fda-8852-3a931081d1e1.md",
    "",
    "",
    "c4e3a62b-d9fe-4fda-8852-3a931081d1e1.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3538982758f4a116054fa89e59de4f03d20fb81b/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/c4e3a62b-d9fe-4fda-8852-3a931081d1e1.05d347036dc41e7d5c32e0cb191cee2934055b46.zh-cn.xlf",
    "",
    "",
    "c4e3a62b-d9fe-4fda-8852-3a931081d1e1.05d347036dc41e7d5c32e0cb191cee2934055b46.zh-cn.xlf"
)
$wsZhCn.Range("H3").Value = "2016-03-22 12:11:49"

# --- de-de sheet: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3d274448feffc1782fe7b51ae4b7508c91d149bb/e2e/8351b26a-afed-4463-8200-ed7236f08299.md",
    "",
    "",
    "8351b26a-afed-4463-8200-ed7236f08299.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1c535adb7e48f7a3b65d108b003e7b122c9cf027/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/8351b26a-afed-4463-8200-ed7236f08299.1f1eec3eddfab6c04388319c08bf8fbb2f1611cc.de-de.xlf",
    "",
    "",
    "8351b26a-afed-4463-8200-ed7236f08299.1f1eec3eddfab6c04388319c08bf8fbb2f1611cc.de-de.xlf"
)
$wsDeDe.Range("H2").Value = "2016-03-22 12:12:03"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3d274448feffc1782fe7b51ae4b7508c91d149bb/e2e/c4e3a62b-d9fe-4fda-8852-3a931081d1e1.md",
    "",
    "",
    "c4e3a62b-d9fe-4fda-8852-3a931081d1e1.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1c535adb7e48f7a3b65d108b003e7b122c9cf027/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/c4e3a62b-d9fe-4fda-8852-3a931081d1e1.05d347036dc41e7d5c32e0cb191cee2934055b46.de-de.xlf",
    "",
    "",
    "c4e3a62b-d9fe-4fda-8852-3a931081d1e1.05d347036dc41e7d5c32e0cb191cee2934055b46.de-de.xlf"
)
$wsDeDe.Range("H3").Value = "2016-03-22 12:12:03"
